$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Toll")

# --- Toll sheet: D21 gets new text "Per Car" ---
$ws2.Range("D21").Value = "Per Car"

# --- Sheet1: D20 formula gains B5 term ---
$ws1.Range("D20").Formula = "=SUM(B3,B4,B5,D9,G17,J8,M11,P11,S11)"

# --- Sheet1: K18 becomes a formula referencing J18 instead of a static value ---
$ws1.Range("K18").Formula = "=J18"
